$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dhw_effic value for Electricity (row 2, column G) from 1 to 0.93
$ws.Range("G2").Value = 0.93

# Update the selected cell/range shown in the sheet view
$ws.Range("G3").Select()
